# #5: fund, bonds, otherbonds, antique done
#
# Sheet "具有相當價值之財產" (property of considerable value, worksheet #5):
#   - row 1 becomes a header row (name/quantity/owner/total/property_category/
#     category/date/legislator_name/legislator_id/source_file/index)
#   - columns F:L are added to every data row (rows 2-18) carrying the
#     normalized metadata (category/date/legislator info/source file/index)
#   - a couple of stray values get corrected
# Sheet "債權" (claims, worksheet #7):
#   - a stray trailing glyph is removed from a date string

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 5 : 具有相當價值之財產
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

# Turn row 1 from a (duplicated) data row into the column-header row.
$ws5.Range("B1").Value = "name"
$ws5.Range("C1").Value = "quantity"
$ws5.Range("D1").Value = "owner"
$ws5.Range("E1").Value = "total"

# New header cells F1:L1 - copy style from an existing header cell first so
# they pick up the same (bold/bordered) formatting, then set their text.
$ws5.Range("B1").Copy($ws5.Range("F1:L1"))
$ws5.Range("F1").Value = "property_category"
$ws5.Range("G1").Value = "category"
$ws5.Range("H1").Value = "date"
$ws5.Range("I1").Value = "legislator_name"
$ws5.Range("J1").Value = "legislator_id"
$ws5.Range("K1").Value = "source_file"
$ws5.Range("L1").Value = "index"

# Fix stray "■" glyph appended to a name.
$ws5.Range("B3").Value = "鑽耳環"

# Corrected quantities.
$ws5.Range("C5").Value = 1
$ws5.Range("C6").Value = 1

# Quantity becomes a textual "s" instead of the numeric 1.
$ws5.Range("C10").Value = "s"

# Populate the new F:L columns for every data row (2-18), using style copied
# from the already-bordered column E of the same row, then overwrite values.
for ($r = 2; $r -le 18; $r++) {
    $ws5.Range("E$r").Copy($ws5.Range("F$r`:L$r"))
    $ws5.Range("F$r").Value = "otherbonds"
    $ws5.Range("G$r").Value = "normal"
    $ws5.Range("H$r").Value = "2013-12-12"
    $ws5.Range("I$r").Value = "張慶忠"
    $ws5.Range("J$r").Value = 1347
    $ws5.Range("K$r").Value = "tmpe4561"
    $indexVal = $ws5.Range("A$r").Value()
    $ws5.Range("L$r").Value = $indexVal
}

Write-Host "sheet5 done"

# ---------------------------------------------------------------------------
# Sheet 7 : 債權
# ---------------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item(7)

# Fix stray "■" glyph appended to a date.
$ws7.Range("F3").Value = "100年08月15日"

Write-Host "sheet7 done"
